$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates used in this batch of "date started" / "date finished" entries.
$d1107 = Get-Date -Year 2018 -Month 11 -Day 7 -Hour 0 -Minute 0 -Second 0
$d1106 = Get-Date -Year 2018 -Month 11 -Day 6 -Hour 0 -Minute 0 -Second 0

# Rows 72-108: species reviewed on 2018-11-07 with 0 results making it into the dataset.
$rows0 = 72..108
foreach ($r in $rows0) {
    $dCell = $ws.Cells.Item($r, 4)   # column D - date started
    $dCell.NumberFormat = "d-mmm"
    $dCell.Value = $d1107

    $eCell = $ws.Cells.Item($r, 5)   # column E - date finished
    $eCell.NumberFormat = "d-mmm"
    $eCell.Value = $d1107

    $gCell = $ws.Cells.Item($r, 7)   # column G - number included in dataset
    $gCell.Value = 0
}

# Row 109 is the odd one out: started 2018-11-06, finished 2018-11-07, 1 made it in.
$ws.Cells.Item(109, 4).NumberFormat = "d-mmm"
$ws.Cells.Item(109, 4).Value = $d1106
$ws.Cells.Item(109, 5).NumberFormat = "d-mmm"
$ws.Cells.Item(109, 5).Value = $d1107
$ws.Cells.Item(109, 7).Value = 1

# Reflect where the reviewer had scrolled/selected when the workbook was saved.
$ws.Range("A71").Select()
